$wb = $excel.ActiveWorkbook

# ---- Sheet "518" ----------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws518 = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws518.Name = "518"

$ws518.Range("B1").Value = "Mo Implanted Target"
$ws518.Range("C1").Value = "Ta Implanted Target"
$ws518.Range("D1").Value = "W Implanted Target"
$ws518.Range("E1").Value = "Weighted Average"

$ws518.Range("A2").Value = '$F(\tau)_{5.18}$'
$ws518.Range("B2").Value = '$9.04 \pm 0.013$'
$ws518.Range("C2").Value = '$0.911 \pm 0.016$'
$ws518.Range("D2").Value = '$9.12 \pm 0.015$'

$bom = [char]0xFEFF
$moTau518 = '$' + $bom + '7.1^{+4.8}_{-2.3}$'

$ws518.Range("A3").Value = '$\tau_{5.18}$ (fs)'
$ws518.Range("B3").Value = $moTau518
$ws518.Range("C3").Value = '$7.1 \pm 5.5$'
$ws518.Range("D3").Value = '$8.0 \pm 6.7$'
$ws518.Range("E3").Value = '$7.5 \pm 3.0$'

$ws518.Rows.Item(3).RowHeight = 17

$f518 = $ws518.Range("C3:E3").Font
$f518.Name = "Arial Unicode MS"
$f518.Size = 10
$f518.Color = 3984993

# ---- Sheet "617" ----------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws617 = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws617.Name = "617"

$ws617.Range("B1").Value = "Mo Implanted Target"
$ws617.Range("C1").Value = "Ta Implanted Target"
$ws617.Range("D1").Value = "W Implanted Target"
$ws617.Range("E1").Value = "Weighted Average"

$ws617.Range("A2").Value = '$F(\tau)_{6.17}$'
$ws617.Range("B2").Value = '$0.992 \pm 0.014$'
$ws617.Range("C2").Value = '$0.976 \pm 0.017$'
$ws617.Range("D2").Value = '$0.988 \pm 0.016$'

$ws617.Range("A3").Value = '$\tau_{6.17}$ (fs)'
$ws617.Range("B3").Value = '$0.4^{+0.7}_{-0.4}$'
$ws617.Range("C3").Value = '$1.4 \pm 1.0$'
$ws617.Range("D3").Value = '$0.6^{+0.9}_{-0.6}$'
$ws617.Range("E3").Value = '$0.7 \pm 0.5$'

$ws617.Rows.Item(3).RowHeight = 17

$f617 = $ws617.Range("C3:E3").Font
$f617.Name = "Arial Unicode MS"
$f617.Size = 10
$f617.Color = 3984993

# ---- Sheet "679" ----------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws679 = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws679.Name = "679"

$ws679.Range("B1").Value = "Mo Implanted Target"
$ws679.Range("C1").Value = "Ta Implanted Target"
$ws679.Range("D1").Value = "W Implanted Target"
$ws679.Range("E1").Value = "Weighted Average"

$ws679.Range("A2").Value = '$F(\tau)_{6.79}$'
$ws679.Range("B2").Value = '$0.995 \pm 0.019$'
$ws679.Range("C2").Value = '$ 0.983 \pm 0.019$'
$ws679.Range("D2").Value = '$0.978 \pm 0.015$'

$ws679.Range("A3").Value = '$\tau_{6.79}$ (fs)'
$ws679.Range("B3").Value = '$0.2^{+0.7}_{-0.2}$'
$ws679.Range("C3").Value = '$0.7^{+0.9}_{-0.7}$'
$ws679.Range("D3").Value = '$0.9 \pm 0.6$'
$ws679.Range("E3").Value = '$0.6 \pm 0.4$'

$ws679.Rows.Item(3).RowHeight = 17

$f679 = $ws679.Range("C3:E3").Font
$f679.Name = "Arial Unicode MS"
$f679.Size = 10
$f679.Color = 3984993

# ---- Selections / active sheet --------------------------------------------
[void]$ws679.Range("D4").Select()
[void]$ws617.Activate()
[void]$ws617.Range("A3").Select()
[void]$ws518.Activate()
[void]$ws518.Range("A3").Select()
